# Applies the edits described by the diff:
#  - Sheet2: insert a new task row "Gör powerpointpresentationen" (status
#    "Påbörjad", 7h actual time) above the "Summa" summary block, update the
#    "Ändra grundfärgerna..." row status to "Påbörjad" with 4h actual time,
#    and bump the "Summa" total (actual time) to 34.
#  - Update remembered cell selections on Sheet2, "Iteration 2" and
#    "Iteration 6" sheets.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")

# Insert a new row before row 24 (the "Summa" row), shifting the existing
# summary rows (Summa / Tid sedan föregående iteration / Total projekttid)
# down by one, the same way Excel's "Insert Row" command would.
$ws2.Rows.Item(24).Insert()

# Update the "Ändra grundfärgerna på knappar och bakgrund." row (row 20):
# status changes from "Ej påbörjad" to "Påbörjad", and actual time of 4h is
# logged.
$ws2.Range("C20").Value = "Påbörjad"
$ws2.Range("E20").Value = 4

# Fill in the newly inserted row 23 with the new task.
$ws2.Range("B23").Value = "Gör powerpointpresentationen"
$ws2.Range("C23").Value = "Påbörjad"
$ws2.Range("E23").Value = 7

# The "Summa" row (now row 25) gets an updated actual-time total.
$ws2.Range("E25").Value = 34

# Restore the remembered selections for the affected sheets.
$ws2.Range("D21").Select()

$wsIt2 = $wb.Worksheets.Item("Iteration 2")
$wsIt2.Range("B32").Select()

$wsIt6 = $wb.Worksheets.Item("Iteration 6")
$wsIt6.Range("C31").Select()

$ws2.Activate()
$ws2.Range("D21").Select()
